$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 2590.1904
$ws.Range("I40").Value = 2568.7693
$ws.Range("J40").Value = 2625
$ws.Range("K40").Value = 2568.7693
$ws.Range("L40").Value = 2625
$ws.Range("M40").Value = -2393.7693
$ws.Range("N40").Value = -2975
$ws.Range("H98").Value = 500
$ws.Range("I98").Value = 500
$ws.Range("K98").Value = 500
$ws.Range("M98").Value = 998
$ws.Range("H122").Value = 500
$ws.Range("I122").Value = 500
$ws.Range("K122").Value = 1500
$ws.Range("M122").Value = 950
# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5289.2563
$ws.Range("I32").Value = 5438.3887
$ws.Range("K32").Value = 5438.3887
$ws.Range("M32").Value = -5151.3887
$ws.Range("H74").Value = 297529.5
$ws.Range("I74").Value = 297529.5
$ws.Range("K74").Value = 297529.5
$ws.Range("M74").Value = -296655.5
$ws.Range("H77").Value = 297529.5
$ws.Range("I77").Value = 297529.5
$ws.Range("K77").Value = 1487647.5
$ws.Range("M77").Value = -1483279.5
$ws.Range("H122").Value = 3540.5334
$ws.Range("I122").Value = 2856.889
$ws.Range("K122").Value = 8570.667000000001
$ws.Range("M122").Value = -6120.667000000001
$ws.Range("H132").Value = 71956.375
$ws.Range("I132").Value = 14656.954
$ws.Range("J132").Value = 702250
$ws.Range("K132").Value = 43970.862
$ws.Range("L132").Value = 2106750
$ws.Range("M132").Value = -41440.862
$ws.Range("N132").Value = -2111810
# --- BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 20185.438
$ws.Range("I86").Value = 9714
$ws.Range("K86").Value = 9714
$ws.Range("M86").Value = -8591
$ws.Range("H89").Value = 20185.438
$ws.Range("I89").Value = 9714
$ws.Range("K89").Value = 48570
$ws.Range("M89").Value = -42954
$ws.Range("H94").Value = 1315.6666
$ws.Range("I94").Value = 667.087
$ws.Range("J94").Value = 2463.1538
$ws.Range("K94").Value = 667.087
$ws.Range("L94").Value = 2463.1538
$ws.Range("M94").Value = -216.087
$ws.Range("N94").Value = -3365.1538
$ws.Range("H95").Value = 18086.7
$ws.Range("J95").Value = 18086.7
$ws.Range("L95").Value = 18086.7
$ws.Range("N95").Value = -23578.7
$ws.Range("H105").Value = 8058.854
$ws.Range("I105").Value = 14431.6875
$ws.Range("J105").Value = 4872.4375
$ws.Range("K105").Value = 14431.6875
$ws.Range("L105").Value = 4872.4375
$ws.Range("M105").Value = -12684.6875
$ws.Range("N105").Value = -8366.4375
$ws.Range("H107").Value = 1585.65
$ws.Range("I107").Value = 1481.138
$ws.Range("J107").Value = 1861.1818
$ws.Range("K107").Value = 1481.138
$ws.Range("L107").Value = 1861.1818
$ws.Range("M107").Value = 438.8620000000001
$ws.Range("N107").Value = -5701.1818
$ws.Range("H134").Value = 3661
$ws.Range("I134").Value = 985.3333
$ws.Range("K134").Value = 2955.9999
$ws.Range("M134").Value = -420.9998999999998
# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3831.111
$ws.Range("I31").Value = 2270.3333
$ws.Range("J31").Value = 5391.8887
$ws.Range("K31").Value = 2270.3333
$ws.Range("L31").Value = 5391.8887
$ws.Range("M31").Value = -1975.3333
$ws.Range("N31").Value = -5981.8887
$ws.Range("H34").Value = 3831.111
$ws.Range("I34").Value = 2270.3333
$ws.Range("J34").Value = 5391.8887
$ws.Range("K34").Value = 2270.3333
$ws.Range("L34").Value = 5391.8887
$ws.Range("M34").Value = -2068.3333
$ws.Range("N34").Value = -5795.8887
$ws.Range("H99").Value = 4438.0586
$ws.Range("I99").Value = 5378.4
$ws.Range("J99").Value = 3094.7144
$ws.Range("K99").Value = 5378.4
$ws.Range("L99").Value = 3094.7144
$ws.Range("M99").Value = -3880.4
$ws.Range("N99").Value = -6090.7144
$ws.Range("H105").Value = 14191
$ws.Range("I105").Value = 4487.375
$ws.Range("J105").Value = 53005.5
$ws.Range("K105").Value = 4487.375
$ws.Range("L105").Value = 53005.5
$ws.Range("M105").Value = -2740.375
$ws.Range("N105").Value = -56499.5
$ws.Range("H107").Value = 442.1111
$ws.Range("I107").Value = 410
$ws.Range("J107").Value = 699
$ws.Range("K107").Value = 410
$ws.Range("L107").Value = 699
$ws.Range("M107").Value = 1510
$ws.Range("N107").Value = -4539
$ws.Range("H126").Value = 4438.0586
$ws.Range("I126").Value = 5378.4
$ws.Range("J126").Value = 3094.7144
$ws.Range("K126").Value = 16135.2
$ws.Range("L126").Value = 9284.143199999999
$ws.Range("M126").Value = -13665.2
$ws.Range("N126").Value = -14224.1432
# --- CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 2757.5
$ws.Range("I5").Value = 1112
$ws.Range("K5").Value = 3336
$ws.Range("M5").Value = -3224
$ws.Range("H39").Value = 2394.6
$ws.Range("I39").Value = 1416.6666
$ws.Range("J39").Value = 2813.7144
$ws.Range("K39").Value = 4249.9998
$ws.Range("L39").Value = 8441.143199999999
$ws.Range("M39").Value = -3955.9998
$ws.Range("N39").Value = -9029.143199999999
$ws.Range("H55").Value = 10809666
$ws.Range("I55").Value = 1800159.8
$ws.Range("K55").Value = 5400479.4
$ws.Range("M55").Value = -5400302.4
$ws.Range("H68").Value = 0
$ws.Range("I68").Value = 0
$ws.Range("K68").Value = 0
$ws.Range("M68").ClearContents()
$ws.Range("H71").Value = 0
$ws.Range("I71").Value = 0
$ws.Range("K71").Value = 0
$ws.Range("M71").ClearContents()
$ws.Range("H113").Value = 942.3158
$ws.Range("I113").Value = 398.57144
$ws.Range("J113").Value = 2464.8
$ws.Range("K113").Value = 1195.71432
$ws.Range("L113").Value = 7394.400000000001
$ws.Range("M113").Value = 974.28568
$ws.Range("N113").Value = -11734.4
$ws.Range("H135").Value = 2757.5
$ws.Range("I135").Value = 1112
$ws.Range("K135").Value = 10008
$ws.Range("M135").Value = -7473
$ws.Range("H139").Value = 3060
$ws.Range("I139").Value = 3060
$ws.Range("J139").Value = 0
$ws.Range("K139").Value = 9180
$ws.Range("L139").Value = 0
$ws.Range("M139").Value = -4040
$ws.Range("N139").ClearContents()
# --- GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H94").Value = 26142.428
$ws.Range("J94").Value = 26142.428
$ws.Range("L94").Value = 26142.428
$ws.Range("N94").Value = -27494.428
$ws.Range("H102").Value = 4775.9585
$ws.Range("I102").Value = 2542.8
$ws.Range("K102").Value = 2542.8
$ws.Range("M102").Value = -920.8000000000002
$ws.Range("H122").Value = 1284.7273
$ws.Range("I122").Value = 1084
$ws.Range("K122").Value = 3252
$ws.Range("M122").Value = -802
$ws.Range("H132").Value = 1527.2778
$ws.Range("I132").Value = 1523.8125
$ws.Range("J132").Value = 1555
$ws.Range("K132").Value = 4571.4375
$ws.Range("L132").Value = 4665
$ws.Range("M132").Value = -2041.4375
$ws.Range("N132").Value = -9725
# --- LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 3448.5
$ws.Range("I40").Value = 1897
$ws.Range("K40").Value = 1897
$ws.Range("M40").Value = -1761
$ws.Range("H123").Value = 69607.55499999999
$ws.Range("J123").Value = 74996
$ws.Range("L123").Value = 74996
$ws.Range("N123").Value = -84796
# --- WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H8").Value = 0
$ws.Range("J8").Value = 0
$ws.Range("L8").Value = 0
$ws.Range("N8").ClearContents()
$ws.Range("H26").Value = 49924.75
$ws.Range("I26").Value = 49900
$ws.Range("K26").Value = 49900
$ws.Range("M26").Value = -49607
$ws.Range("H109").Value = 92666.664
$ws.Range("I109").Value = 89000
$ws.Range("J109").Value = 100000
$ws.Range("K109").Value = 89000
$ws.Range("L109").Value = 100000
$ws.Range("M109").Value = -87613
$ws.Range("N109").Value = -102774
$ws.Range("H113").Value = 1085.7084
$ws.Range("I113").Value = 978.0526
$ws.Range("K113").Value = 2934.1578
$ws.Range("M113").Value = -764.1578
$ws.Range("H132").Value = 2359.8
$ws.Range("I132").Value = 2674.125
$ws.Range("J132").Value = 2000.5714
$ws.Range("K132").Value = 8022.375
$ws.Range("L132").Value = 6001.7142
$ws.Range("M132").Value = -5492.375
$ws.Range("N132").Value = -11061.7142
$ws.Range("H136").Value = 5120.727
$ws.Range("I136").Value = 3310.7334
$ws.Range("K136").Value = 9932.200199999999
$ws.Range("M136").Value = -7382.200199999999
